# Update the cryptos worksheet with freshly scraped price/volume data.
# For cells whose new value looks like a plain number (e.g. "233.28"),
# we force a Text number format before assignment and reset the style
# back to Normal afterwards so the cell keeps its original (default)
# style while the value is preserved verbatim as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '92.193.13'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '3.099.68'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  +0.05%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '233.28'
$c.Style = 'Normal'
$ws.Range("E5").Value = '  -2.90%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '612.24'
$c.Style = 'Normal'
$ws.Range("E6").Value = '  -1.24%  '
$ws.Range("E7").Value = '  -3.39%  '
$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '0.385'
$c.Style = 'Normal'
$ws.Range("E8").Value = '  -0.45%  '
$ws.Range("D10").Value = '3.095.76'
$ws.Range("E10").Value = '  -1.92%  '
$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '0.774'
$c.Style = 'Normal'
$ws.Range("E11").Value = '  +4.17%  '
$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '0.197'
$c.Style = 'Normal'
$ws.Range("E12").Value = '  -3.60%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.0000243'
$c.Style = 'Normal'
$ws.Range("E13").Value = '  -4.11%  '
$ws.Range("D14").Value = '91.974.55'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '33.69'
$c.Style = 'Normal'
$ws.Range("E15").Value = '  -4.05%  '
$ws.Range("B16").Value = 'Toncoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '5.39'
$c.Style = 'Normal'
$ws.Range("E16").Value = '  -3.64%  '
$ws.Range("D17").Value = '3.674.44'
$ws.Range("E17").Value = '  -1.56%  '
$ws.Range("D18").Value = '3.086.03'
$ws.Range("E18").Value = '  -2.50%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '3.80'
$c.Style = 'Normal'
$ws.Range("E19").Value = '  +1.84%  '
$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '14.37'
$c.Style = 'Normal'
$ws.Range("E20").Value = '  -4.29%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '5.73'
$c.Style = 'Normal'
$ws.Range("E21").Value = '  -3.13%  '
$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '0.0000198'
$c.Style = 'Normal'
$ws.Range("E22").Value = '  -2.87%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '434.59'
$c.Style = 'Normal'
$ws.Range("E23").Value = '  -4.97%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '9.06'
$c.Style = 'Normal'
$ws.Range("E24").Value = '  -1.48%  '
$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '5.55'
$c.Style = 'Normal'
$ws.Range("E25").Value = '  -6.03%  '
$c = $ws.Range("D26")
$c.NumberFormat = '@'
$c.Value = '85.15'
$c.Style = 'Normal'
$ws.Range("E26").Value = '  -4.02%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '11.32'
$c.Style = 'Normal'
$ws.Range("E27").Value = '  -4.13%  '
$ws.Range("D28").Value = '3.259.27'
$ws.Range("E28").Value = '  -1.76%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range("E29").Value = '  +0.03%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '0.178'
$c.Style = 'Normal'
$ws.Range("E30").Value = '  +6.31%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '0.232'
$c.Style = 'Normal'
$ws.Range("E31").Value = '  +1.41%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '0.124'
$c.Style = 'Normal'
$ws.Range("E32").Value = '  -16.26%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '1.04'
$c.Style = 'Normal'
$ws.Range("E33").Value = '  -35.62%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '9.11'
$c.Style = 'Normal'
$ws.Range("E34").Value = '  -2.92%  '
$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '7.94'
$c.Style = 'Normal'
$ws.Range("E35").Value = '  +6.15%  '
$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '0.155'
$c.Style = 'Normal'
$ws.Range("E36").Value = '  -12.08%  '
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '25.46'
$c.Style = 'Normal'
$ws.Range("E37").Value = '  -3.44%  '
$ws.Range("E38").Value = '  +0.19%  '
$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '1.88'
$c.Style = 'Normal'
$ws.Range("E39").Value = '  -3.28%  '
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '23.85'
$c.Style = 'Normal'
$ws.Range("E40").Value = '  +7.60%  '
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '0.434'
$c.Style = 'Normal'
$ws.Range("E41").Value = '  -2.08%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '1.27'
$c.Style = 'Normal'
$ws.Range("E42").Value = '  -4.20%  '
$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '465.15'
$c.Style = 'Normal'
$ws.Range("E43").Value = '  -5.66%  '
$ws.Range("E44").Value = '  -3.36%  '
$ws.Range("E45").Value = '  +0.14%  '
$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '158.58'
$c.Style = 'Normal'
$ws.Range("E46").Value = '  +1.85%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '0.678'
$c.Style = 'Normal'
$ws.Range("E47").Value = '  -4.45%  '
$ws.Range("E48").Value = '  -5.80%  '
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '43.77'
$c.Style = 'Normal'
$ws.Range("E49").Value = '  -0.79%  '
$ws.Range("E50").Value = '  -3.50%  '
$ws.Range("E51").Value = '  -0.06%  '

Write-Host "Applied cryptos list update."
